$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 315.42856
$ws.Cells.Item(18, 9).Value = 315.42856
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 11).Value = 315.42856
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 13).Value = -31.42856
$ws.Cells.Item(18, 14).ClearContents()

$ws.Cells.Item(111, 8).Value = 72630.92999999999
$ws.Cells.Item(111, 9).Value = 92030.45
$ws.Cells.Item(111, 10).Value = 1499.3334
$ws.Cells.Item(111, 11).Value = 276091.35
$ws.Cells.Item(111, 12).Value = 4498.0002
$ws.Cells.Item(111, 13).Value = -273024.35
$ws.Cells.Item(111, 14).Value = -10632.0002

$ws.Cells.Item(113, 8).Value = 1900.375
$ws.Cells.Item(113, 9).Value = 2160.8
$ws.Cells.Item(113, 10).Value = 1466.3334
$ws.Cells.Item(113, 11).Value = 2160.8
$ws.Cells.Item(113, 12).Value = 1466.3334
$ws.Cells.Item(113, 13).Value = 1093.2
$ws.Cells.Item(113, 14).Value = -7974.3334

$ws.Cells.Item(138, 8).Value = 6908.5615
$ws.Cells.Item(138, 9).Value = 4243.909
$ws.Cells.Item(138, 10).Value = 7545.7607
$ws.Cells.Item(138, 11).Value = 12731.727
$ws.Cells.Item(138, 12).Value = 22637.2821
$ws.Cells.Item(138, 13).Value = -7591.726999999999
$ws.Cells.Item(138, 14).Value = -32917.2821

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5904.8335
$ws.Cells.Item(32, 9).Value = 5904.8335
$ws.Cells.Item(32, 11).Value = 5904.8335
$ws.Cells.Item(32, 13).Value = -5617.8335

$ws.Cells.Item(45, 8).Value = 2449.5715
$ws.Cells.Item(45, 9).Value = 2071.1667
$ws.Cells.Item(45, 10).Value = 2733.375
$ws.Cells.Item(45, 11).Value = 2071.1667
$ws.Cells.Item(45, 12).Value = 2733.375
$ws.Cells.Item(45, 13).Value = -1694.1667
$ws.Cells.Item(45, 14).Value = -3487.375

$ws.Cells.Item(61, 8).Value = 3393.2632
$ws.Cells.Item(61, 9).Value = 3563.4119
$ws.Cells.Item(61, 10).Value = 1947
$ws.Cells.Item(61, 11).Value = 3563.4119
$ws.Cells.Item(61, 12).Value = 1947
$ws.Cells.Item(61, 13).Value = -3351.4119
$ws.Cells.Item(61, 14).Value = -2371

$ws.Cells.Item(132, 8).Value = 5225.4688
$ws.Cells.Item(132, 9).Value = 4844.1304
$ws.Cells.Item(132, 11).Value = 14532.3912
$ws.Cells.Item(132, 13).Value = -12002.3912

$ws.Cells.Item(136, 8).Value = 3393.2632
$ws.Cells.Item(136, 9).Value = 3563.4119
$ws.Cells.Item(136, 10).Value = 1947
$ws.Cells.Item(136, 11).Value = 10690.2357
$ws.Cells.Item(136, 12).Value = 5841
$ws.Cells.Item(136, 13).Value = -8140.235700000001
$ws.Cells.Item(136, 14).Value = -10941

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(24, 8).Value = 19996
$ws.Cells.Item(24, 9).Value = 19996
$ws.Cells.Item(24, 11).Value = 19996
$ws.Cells.Item(24, 13).Value = -19761

$ws.Cells.Item(25, 8).Value = 6882
$ws.Cells.Item(25, 9).Value = 6882
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 11).Value = 6882
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 13).Value = -6647
$ws.Cells.Item(25, 14).ClearContents()

$ws.Cells.Item(86, 8).Value = 2129375
$ws.Cells.Item(86, 9).Value = 2129375
$ws.Cells.Item(86, 11).Value = 2129375
$ws.Cells.Item(86, 13).Value = -2128252

$ws.Cells.Item(89, 8).Value = 2129375
$ws.Cells.Item(89, 9).Value = 2129375
$ws.Cells.Item(89, 11).Value = 10646875
$ws.Cells.Item(89, 13).Value = -10641259

$ws.Cells.Item(134, 8).Value = 35878.625
$ws.Cells.Item(134, 9).Value = 5635.65
$ws.Cells.Item(134, 10).Value = 86283.586
$ws.Cells.Item(134, 11).Value = 16906.95
$ws.Cells.Item(134, 12).Value = 258850.758
$ws.Cells.Item(134, 13).Value = -14371.95
$ws.Cells.Item(134, 14).Value = -263920.758

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(3, 8).Value = 6999.6665
$ws.Cells.Item(3, 9).Value = 1000
$ws.Cells.Item(3, 10).Value = 9999.5
$ws.Cells.Item(3, 11).Value = 1000
$ws.Cells.Item(3, 12).Value = 9999.5
$ws.Cells.Item(3, 13).Value = -887
$ws.Cells.Item(3, 14).Value = -10225.5

$ws.Cells.Item(16, 8).Value = 3069.1724
$ws.Cells.Item(16, 9).Value = 2662.5
$ws.Cells.Item(16, 11).Value = 2662.5
$ws.Cells.Item(16, 13).Value = -2375.5

$ws.Cells.Item(17, 8).Value = 5000
$ws.Cells.Item(17, 10).Value = 5000
$ws.Cells.Item(17, 12).Value = 5000
$ws.Cells.Item(17, 14).Value = -5348

$ws.Cells.Item(25, 8).Value = 5999
$ws.Cells.Item(25, 10).Value = 5999
$ws.Cells.Item(25, 12).Value = 5999
$ws.Cells.Item(25, 14).Value = -6347

$ws.Cells.Item(31, 8).Value = 47758.547
$ws.Cells.Item(31, 9).Value = 1593.8334
$ws.Cells.Item(31, 11).Value = 1593.8334
$ws.Cells.Item(31, 13).Value = -1298.8334

$ws.Cells.Item(34, 8).Value = 47758.547
$ws.Cells.Item(34, 9).Value = 1593.8334
$ws.Cells.Item(34, 11).Value = 1593.8334
$ws.Cells.Item(34, 13).Value = -1391.8334

$ws.Cells.Item(58, 8).Value = 2631.5715
$ws.Cells.Item(58, 9).Value = 2701.8333
$ws.Cells.Item(58, 10).Value = 2210
$ws.Cells.Item(58, 11).Value = 2701.8333
$ws.Cells.Item(58, 12).Value = 2210
$ws.Cells.Item(58, 13).Value = -2498.8333
$ws.Cells.Item(58, 14).Value = -2616

$ws.Cells.Item(113, 8).Value = 3069.1724
$ws.Cells.Item(113, 9).Value = 2662.5
$ws.Cells.Item(113, 11).Value = 2662.5
$ws.Cells.Item(113, 13).Value = -492.5

$ws.Cells.Item(122, 8).Value = 0
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).ClearContents()
$ws.Cells.Item(122, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 1562.6976
$ws.Cells.Item(132, 9).Value = 1356.0541
$ws.Cells.Item(132, 11).Value = 4068.1623
$ws.Cells.Item(132, 13).Value = -1538.1623

$ws.Cells.Item(136, 8).Value = 2631.5715
$ws.Cells.Item(136, 9).Value = 2701.8333
$ws.Cells.Item(136, 10).Value = 2210
$ws.Cells.Item(136, 11).Value = 8105.499899999999
$ws.Cells.Item(136, 12).Value = 6630
$ws.Cells.Item(136, 13).Value = -5555.499899999999
$ws.Cells.Item(136, 14).Value = -11730

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 1792137.4
$ws.Cells.Item(4, 10).Value = 302065
$ws.Cells.Item(4, 12).Value = 906195
$ws.Cells.Item(4, 14).Value = -906419

$ws.Cells.Item(17, 8).Value = 93.666664
$ws.Cells.Item(17, 9).Value = 91
$ws.Cells.Item(17, 10).Value = 99
$ws.Cells.Item(17, 11).Value = 273
$ws.Cells.Item(17, 12).Value = 297
$ws.Cells.Item(17, 13).Value = -104
$ws.Cells.Item(17, 14).Value = -635

$ws.Cells.Item(121, 8).Value = 1002306.3
$ws.Cells.Item(121, 10).Value = 1252070.4
$ws.Cells.Item(121, 12).Value = 3756211.2
$ws.Cells.Item(121, 14).Value = -3758831.2

$ws.Cells.Item(138, 8).Value = 2832.8823
$ws.Cells.Item(138, 9).Value = 2243.1538
$ws.Cells.Item(138, 11).Value = 6729.4614
$ws.Cells.Item(138, 13).Value = -1589.4614

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(12, 8).Value = 4000
$ws.Cells.Item(12, 10).Value = 4000
$ws.Cells.Item(12, 12).Value = 4000
$ws.Cells.Item(12, 14).Value = -4280

$ws.Cells.Item(14, 9).Value = 15412941
$ws.Cells.Item(14, 10).Value = 1722506.6
$ws.Cells.Item(14, 11).Value = 15412941
$ws.Cells.Item(14, 12).Value = 1722506.6
$ws.Cells.Item(14, 13).Value = -15412773
$ws.Cells.Item(14, 14).Value = -1722842.6

$ws.Cells.Item(102, 8).Value = 3384.8
$ws.Cells.Item(102, 9).Value = 2032.1333
$ws.Cells.Item(102, 10).Value = 5413.8
$ws.Cells.Item(102, 11).Value = 2032.1333
$ws.Cells.Item(102, 12).Value = 5413.8
$ws.Cells.Item(102, 13).Value = -410.1333
$ws.Cells.Item(102, 14).Value = -8657.799999999999

$ws.Cells.Item(107, 8).Value = 1952.0588
$ws.Cells.Item(107, 9).Value = 1776
$ws.Cells.Item(107, 11).Value = 1776
$ws.Cells.Item(107, 13).Value = 144

$ws.Cells.Item(132, 8).Value = 99174.09
$ws.Cells.Item(132, 9).Value = 9833.444
$ws.Cells.Item(132, 11).Value = 29500.332
$ws.Cells.Item(132, 13).Value = -26970.332

$ws.Cells.Item(133, 8).Value = 55194.5
$ws.Cells.Item(133, 10).Value = 55194.5
$ws.Cells.Item(133, 12).Value = 55194.5
$ws.Cells.Item(133, 14).Value = -65314.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 356
$ws.Cells.Item(16, 9).Value = 356
$ws.Cells.Item(16, 11).Value = 356
$ws.Cells.Item(16, 13).Value = -186

$ws.Cells.Item(20, 8).Value = 815995.5600000001
$ws.Cells.Item(20, 10).Value = 815995.5600000001
$ws.Cells.Item(20, 12).Value = 815995.5600000001
$ws.Cells.Item(20, 14).Value = -816447.5600000001

$ws.Cells.Item(40, 8).Value = 6838.387
$ws.Cells.Item(40, 9).Value = 6427.6
$ws.Cells.Item(40, 11).Value = 6427.6
$ws.Cells.Item(40, 13).Value = -6291.6

$ws.Cells.Item(46, 8).Value = 2069
$ws.Cells.Item(46, 9).Value = 2206.75
$ws.Cells.Item(46, 11).Value = 2206.75
$ws.Cells.Item(46, 13).Value = -2018.75

$ws.Cells.Item(69, 8).Value = 30163
$ws.Cells.Item(69, 10).Value = 30163
$ws.Cells.Item(69, 12).Value = 30163
$ws.Cells.Item(69, 14).Value = -31785

$ws.Cells.Item(72, 8).Value = 30163
$ws.Cells.Item(72, 10).Value = 30163
$ws.Cells.Item(72, 12).Value = 90489
$ws.Cells.Item(72, 14).Value = -98601

$ws.Cells.Item(123, 8).Value = 82966.2
$ws.Cells.Item(123, 10).Value = 82966.2
$ws.Cells.Item(123, 12).Value = 82966.2
$ws.Cells.Item(123, 14).Value = -92766.2

$ws.Cells.Item(136, 8).Value = 284780.22
$ws.Cells.Item(136, 9).Value = 404943.03
$ws.Cells.Item(136, 11).Value = 1214829.09
$ws.Cells.Item(136, 13).Value = -1212279.09

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 3600
$ws.Cells.Item(14, 9).Value = 3600
$ws.Cells.Item(14, 11).Value = 3600
$ws.Cells.Item(14, 13).Value = -3432

$ws.Cells.Item(16, 8).Value = 20420
$ws.Cells.Item(16, 10).Value = 20420
$ws.Cells.Item(16, 12).Value = 20420
$ws.Cells.Item(16, 14).Value = -21004

$ws.Cells.Item(132, 8).Value = 18440.508
$ws.Cells.Item(132, 9).Value = 3373.7964
$ws.Cells.Item(132, 11).Value = 10121.3892
$ws.Cells.Item(132, 13).Value = -7591.389200000001

$ws.Cells.Item(136, 8).Value = 69484.53
$ws.Cells.Item(136, 9).Value = 16213.267
$ws.Cells.Item(136, 10).Value = 335840.84
$ws.Cells.Item(136, 11).Value = 48639.801
$ws.Cells.Item(136, 12).Value = 1007522.52
$ws.Cells.Item(136, 13).Value = -46089.801
$ws.Cells.Item(136, 14).Value = -1012622.52
